$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section 1 (rows 4-15): Multiple Choice Questions hypothesis test
#   H0/H1 labels now quote 48.0% instead of 45.4%; underlying sample data
#   updated (n=150, observed=101) and p0 = 0.48
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "H0:  Agent accuracy = 48.0%"
$ws.Range("B5").Value = "H1: Agent accuracy > 48.0%"
$ws.Range("D7").Value = 0.48
$ws.Range("D9").Value = 150
$ws.Range("D13").Value = 101

# ---------------------------------------------------------------------------
# Section 2 (rows 24-35): Concept Match NOTA vs. Max Freq hypothesis test
#   H0/H1 labels now quote 67.3% instead of 60.2%; n=150, observed=105
# ---------------------------------------------------------------------------
$ws.Range("B24").Value = "H0:  Agent accuracy = 67.3%"
$ws.Range("B25").Value = "H1: Agent accuracy > 67.3%"
$ws.Range("D27").Value = 0.673
$ws.Range("D29").Value = 150
$ws.Range("D33").Value = 105

# ---------------------------------------------------------------------------
# Section 3 (rows 47-55): Concept Match Not vs. Random test -- only the
#   observed count changes (12 -> 8); p0/n stay the same.
# ---------------------------------------------------------------------------
$ws.Range("D53").Value = 8

# ---------------------------------------------------------------------------
# Row 57 conclusion label: this block's result flips to "accept H0"
# ---------------------------------------------------------------------------
$ws.Range("B57").Value = "Conclusion: Accept H0."

# ---------------------------------------------------------------------------
# Section 4 (rows 64-76): Concept Match V2 vs. Max Frequency test
#   H0/H1 labels now quote 56.5% instead of 55%; n=162, observed=106
# ---------------------------------------------------------------------------
$ws.Range("B64").Value = "H0:  Agent accuracy = 56.5%"
$ws.Range("B65").Value = "H1: Agent accuracy > 56.5%"
$ws.Range("D67").Value = 0.565
$ws.Range("D69").Value = 162
$ws.Range("D73").Value = 106

# ---------------------------------------------------------------------------
# Section 5 (rows 87-98): Concept Match V3 vs. Concept Match V2 test
#   H0/H1 labels now quote 47% instead of 55%; n=867, observed=502, and a
#   new helper cell I97 cross-checks the mean using a 0.579 rate.
# ---------------------------------------------------------------------------
$ws.Range("B87").Value = "H0:  Agent accuracy = 47%"
$ws.Range("B88").Value = "H1: Agent accuracy > 47%"
$ws.Range("D92").Value = 867
$ws.Range("D96").Value = 502
$ws.Range("I97").Formula = "=D92*0.579"

# ---------------------------------------------------------------------------
# View state: scroll position / active selection
# ---------------------------------------------------------------------------
$ws.Range("K95").Select()
